$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top and shift everything else down.
$ws.Rows.Item(1).Insert()

# Put the new title string into A1.
$ws.Range("A1").Value = "Writing Intensive WINT Courses Offered in 2016-2017"

# Give the new title cell its own (non-default) font color: plain black,
# matching the new font/cellXf pair added to styles.xml.
$ws.Range("A1").Font.Color = 0

# Reset the view: scroll back to the top and select A5 (matches the
# post-edit sheetView/selection in the diff).
$ws.Application.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("A5").Select() | Out-Null
